$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Femacal de La Calera / Ciboulette was added.
# It belongs chronologically right after the current row 164, so insert a
# fresh row at 165 (pushing the former rows 165-265 down to 166-266) and
# populate it with the new observation.
$ws.Rows.Item(165).Insert()

$ws.Cells.Item(165, 1).Value = 3
$ws.Cells.Item(165, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(165, 3).Value = "Coquimbo"
$ws.Cells.Item(165, 4).Value = 44606
$ws.Cells.Item(165, 5).Value = 5
$ws.Cells.Item(165, 6).Value = 100112039
$ws.Cells.Item(165, 7).Value = "Ciboulette"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 160
$ws.Cells.Item(165, 11).Value = 1500
$ws.Cells.Item(165, 12).Value = 1500
$ws.Cells.Item(165, 13).Value = 1500
$ws.Cells.Item(165, 14).Value = '$/docena de atados'
$ws.Cells.Item(165, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(165, 16).Value = 500
$ws.Cells.Item(165, 17).Value = 3
$ws.Cells.Item(165, 18).Value = "Hortaliza"
